$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, preventing Excel from
# auto-converting numeric-looking strings (e.g. "218.91", "65.30")
# into numbers, and without leaving a residual number-format style
# applied to the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '26.419.15'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.692.18'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  +0.86%  '
Set-TextValue $ws.Range('D5') '218.91'
$ws.Range('E5').Value = '  +1.33%  '
Set-TextValue $ws.Range('D6') '0.5480'
$ws.Range('E6').Value = '  +7.50%  '
$ws.Range('E7').Value = '  +0.86%  '
Set-TextValue $ws.Range('D8') '0.2716'
$ws.Range('E8').Value = '  +1.10%  '
Set-TextValue $ws.Range('D9') '0.06461'
$ws.Range('E9').Value = '  +1.40%  '
Set-TextValue $ws.Range('D10') '22.06'
$ws.Range('E10').Value = '  +1.00%  '
Set-TextValue $ws.Range('D11') '0.07702'
$ws.Range('E11').Value = '  +3.56%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.695.82'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D13') '4.540'
$ws.Range('E13').Value = '  +0.61%  '
Set-TextValue $ws.Range('D14') '0.5821'
$ws.Range('E14').Value = '  +0.20%  '
Set-TextValue $ws.Range('D15') '0.000008384'
$ws.Range('E15').Value = '  -1.28%  '
Set-TextValue $ws.Range('D16') '65.30'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('D17').Value = '26.509.06'
$ws.Range('E17').Value = '  +2.39%  '
Set-TextValue $ws.Range('D18') '4.951'
$ws.Range('E18').Value = '  +0.49%  '
Set-TextValue $ws.Range('D21') '190.04'
$ws.Range('E21').Value = '  +0.28%  '
Set-TextValue $ws.Range('D22') '6.230'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('E23').Value = '  +0.84%  '
Set-TextValue $ws.Range('D24') '150.36'
$ws.Range('E24').Value = '  +3.99%  '
Set-TextValue $ws.Range('D25') '0.1305'
$ws.Range('E25').Value = '  +6.89%  '
Set-TextValue $ws.Range('D26') '7.881'
$ws.Range('E26').Value = '  +3.66%  '
Set-TextValue $ws.Range('D27') '15.72'
$ws.Range('E27').Value = '  +0.29%  '
$ws.Range('E28').Value = '  +6.09%  '
Set-TextValue $ws.Range('D29') '0.06317'
$ws.Range('E29').Value = '  -4.81%  '
Set-TextValue $ws.Range('D30') '1.329'
Set-TextValue $ws.Range('D31') '3.580'
$ws.Range('E31').Value = '  +0.33%  '
Set-TextValue $ws.Range('D32') '3.582'
$ws.Range('E32').Value = '  +1.64%  '
Set-TextValue $ws.Range('D33') '1.674'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('E34').Value = '  +2.39%  '
Set-TextValue $ws.Range('D35') '0.6216'
$ws.Range('E35').Value = '  +0.93%  '
Set-TextValue $ws.Range('D36') '2.417'
$ws.Range('E36').Value = '  +2.07%  '
Set-TextValue $ws.Range('D37') '2.730'
$ws.Range('E37').Value = '  +1.65%  '
Set-TextValue $ws.Range('D38') '6.212'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = '1.120.12'
$ws.Range('E39').Value = '  +2.27%  '
Set-TextValue $ws.Range('D40') '0.01638'
$ws.Range('E40').Value = '  +2.64%  '
Set-TextValue $ws.Range('D41') '0.8798'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('E42').Value = '  +0.71%  '
Set-TextValue $ws.Range('D43') '101.24'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '1.843.60'
$ws.Range('E44').Value = '  +1.55%  '
Set-TextValue $ws.Range('D45') '0.00000000110'
$ws.Range('E45').Value = '  -4.34%  '
Set-TextValue $ws.Range('D46') '57.33'
$ws.Range('E46').Value = '  +1.70%  '
Set-TextValue $ws.Range('D47') '8.202'
$ws.Range('E47').Value = '  +0.99%  '
Set-TextValue $ws.Range('D48') '1.009'
$ws.Range('E48').Value = '  +0.57%  '
Set-TextValue $ws.Range('D49') '0.05277'
$ws.Range('E49').Value = '  +0.78%  '
Set-TextValue $ws.Range('D50') '0.4307'
$ws.Range('E50').Value = '  +0.65%  '
Set-TextValue $ws.Range('D51') '6.062'
$ws.Range('E51').Value = '  +1.18%  '
